$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Proximity sheet: append 3 new "door" event rows (rows 19-21)
# ---------------------------------------------------------------------------
$proximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "18:14:52", "18:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "18:14:58", "18:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "18:15:03", "18:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)

$startRow = 19
for ($i = 0; $i -lt $proximityRows.Length; $i++) {
    $r = $startRow + $i
    $data = $proximityRows[$i]
    # Prefix the date string with an apostrophe so Excel stores it as literal
    # text instead of auto-converting it to a date serial number.
    $proximity.Cells.Item($r, 1).Value = "'" + $data[0]
    $proximity.Cells.Item($r, 2).Value = $data[1]
    $proximity.Cells.Item($r, 3).Value = $data[2]
    $proximity.Cells.Item($r, 4).Value = $data[3]
    $proximity.Cells.Item($r, 5).Value = $data[4]
    $proximity.Cells.Item($r, 6).Value = $data[5]
}

# ---------------------------------------------------------------------------
# Camera sheet: append 3 new "Image Captured" rows (rows 18-20)
# ---------------------------------------------------------------------------
$camera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-02-01", "18:14:53", "18:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "18:14:56", "18:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "18:15:04", "18:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = 18
for ($i = 0; $i -lt $cameraRows.Length; $i++) {
    $r = $startRow + $i
    $data = $cameraRows[$i]
    $camera.Cells.Item($r, 1).Value = "'" + $data[0]
    $camera.Cells.Item($r, 2).Value = $data[1]
    $camera.Cells.Item($r, 3).Value = $data[2]
    $camera.Cells.Item($r, 4).Value = $data[3]
    $camera.Cells.Item($r, 5).Value = $data[4]
    $camera.Cells.Item($r, 6).Value = $data[5]
}
